# Auto-generated edit script: updates crypto price/volume data per the commit diff.
# Uses NumberFormat="@" + ClearFormats() around numeric-looking literals so Excel's
# COM input parser doesn't silently coerce text like "0.999" or "1.00" into numbers
# (which would drop the meaningful trailing zeros / shared-string typing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.514.79"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.639.03"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "524.64"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.65%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "153.43"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.72%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.19%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.576"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.26%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.ClearFormats()
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("E10").Value = "  +3.03%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.347"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "3.103.10"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "60.497.41"
$ws.Range("E14").Value = "  -0.45%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.78"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.10%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000141"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.656.44"
$ws.Range("E17").Value = "  +0.94%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.85%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "351.34"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.74%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.56"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.86%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.25"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.45%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.44%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "61.13"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.34%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.426"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.18%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.166"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.47%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "0.0₃0844"
$ws.Range("E27").Value = "  -0.31%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.ClearFormats()
$ws.Range("E28").Value = "  -2.43%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E29").Value = "  +0.05%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.ClearFormats()
$ws.Range("E30").Value = "  +3.60%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "19.29"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.ClearFormats()
$ws.Range("E32").Value = "  +1.33%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "150.03"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.95%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.06"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  -1.32%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.898"
$c.ClearFormats()
$ws.Range("E36").Value = "  +5.11%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.883"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.53%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "36.74"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.80%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "305.14"
$c.ClearFormats()
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("E40").Value = "  -2.22%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.73"
$c.ClearFormats()
$ws.Range("E41").Value = "  -0.61%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.634"
$c.ClearFormats()
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "19.97"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.61%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0556"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  +1.02%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.82"
$c.ClearFormats()
$ws.Range("E48").Value = "  -3.30%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "10.36"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.49%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.96"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "1.979.65"
$ws.Range("E51").Value = "  -0.94%  "
